# OC0303 tjekBesked.docx
#
# Commit message: "Rettede alle UML filer / Har kun rettede navne til at
# passe, der er stadig aendringer som ikke er blevet aendret."
# ("Fixed all UML files. Only fixed names/wording to match, there are
# still changes that haven't been made yet.")
#
# => Two wording fixes in the body text:
#   1. "UC03 Send Besked"  -> "UC03 Send Beskeder"
#   2. "... mellem klient k og bruger b ..." -> "... mellem patient p og bruger b ..."
#
# (The raw OOXML also shows the built-in style catalogue's internal
# w:styleId values re-minted into Danish spellings, e.g. Heading1 ->
# Overskrift1 / Title -> Titel - with the English w:name left untouched.
# That is what happens when the same content is subsequently re-saved by
# a Danish-language copy of Word; it is not something reachable from the
# Word object model - Style identity/NameLocal can't be used to repoint a
# style's id - so it is left alone here.)

$d = $word.ActiveDocument

# 1) "Send Besked" -> "Send Beskeder"
#    Append "er" right after the existing "Besked" word instead of doing a
#    blind text replace, so the run layout / proofErr bookkeeping around
#    that word stays intact.
$r = $d.Content
$r.Find.Execute("Send Besked", $true, $true, $false, $false, $false,
                 $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter("er")

# 2) "klient k" -> "patient p"
$d.Content.Find.Execute("klient k", $true, $true, $false, $false, $false,
                         $true, 1, $false, "patient p", 2) | Out-Null
